# Update "想去人数" (want-to-go count) values in column F across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1297
$ws.Range("F4").Value = 19
$ws.Range("F5").Value = 7477
$ws.Range("F6").Value = 1876
$ws.Range("F7").Value = 6453
$ws.Range("F8").Value = 151
$ws.Range("F9").Value = 2010
$ws.Range("F10").Value = 545
$ws.Range("F13").Value = 44
$ws.Range("F16").Value = 60
$ws.Range("F17").Value = 8320
$ws.Range("F18").Value = 151
$ws.Range("F19").Value = 61
$ws.Range("F20").Value = 193
$ws.Range("F22").Value = 1789
$ws.Range("F24").Value = 11
$ws.Range("F28").Value = 181
$ws.Range("F30").Value = 1961
$ws.Range("F31").Value = 835
$ws.Range("F32").Value = 450
$ws.Range("F35").Value = 150
$ws.Range("F36").Value = 129
$ws.Range("F38").Value = 3946

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 704
$ws.Range("F5").Value = 297

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1297
$ws.Range("F6").Value = 7477
$ws.Range("F8").Value = 1876
$ws.Range("F9").Value = 6453
$ws.Range("F10").Value = 2010
$ws.Range("F13").Value = 545
$ws.Range("F18").Value = 44
$ws.Range("F22").Value = 60
$ws.Range("F23").Value = 8320
$ws.Range("F24").Value = 61
$ws.Range("F25").Value = 193
$ws.Range("F27").Value = 1789
$ws.Range("F28").Value = 11
$ws.Range("F31").Value = 181
$ws.Range("F32").Value = 1961
$ws.Range("F33").Value = 835
$ws.Range("F35").Value = 450
$ws.Range("F40").Value = 150
$ws.Range("F41").Value = 129

